$d = $word.ActiveDocument

# --- Step 1: plain text substitutions (month names, and "size" -> "six") ---
# "Three to size divers" -> "Three to six divers" is done by removing the
# "ze" from "size" (leaving "si") and later inserting a single "x" run, so
# that only the truly-changed character becomes its own run, matching how
# Word records an in-place word edit.

$rng = $d.Content
$rng.Find.Execute("May", $true, $false, $false, $false, $false, $true, 1, $false, "January", 2)

$rng = $d.Content
$rng.Find.Execute("August", $true, $false, $false, $false, $false, $true, 1, $false, "September", 2)

$rng = $d.Content
$rng.Find.Execute("ze divers", $true, $false, $false, $false, $false, $true, 1, $false, " divers", 2)

# --- Step 2: insert the new "x" right after "si" (completing size -> six) ---
$rng = $d.Content
$rng.Find.Execute("si divers", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(1)
$rng.MoveEnd(1, 2)
$rng.Collapse(0)
$rng.InsertAfter("x")
# Pin this run so the save step does not re-coalesce it into its neighbors.
$rng.Bold = 1
$rng.Bold = 0

# --- Step 3: pin the two inserted month names into their own runs ---
# (Must happen after all other text mutations in this paragraph, since a
# later Find/replace elsewhere in the paragraph would otherwise re-merge an
# already-pinned run back into its neighbor.)
$rng = $d.Content
$rng.Find.Execute("September", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Bold = 1
$rng.Bold = 0

$rng = $d.Content
$rng.Find.Execute("January", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Bold = 1
$rng.Bold = 0
